# "Roll back changes on Visium": remove the preparation_instrument_vendor /
# preparation_instrument_model lookup sheets (and the Visium columns Q/R that
# depended on them), and restore the previous pav:createdOn timestamp on the
# .metadata sheet.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Keep comment authorship consistent with the workbook's existing
# "CEDAR Metadata Validator" comment author.
$excel.UserName = "CEDAR Metadata Validator"

$ws = $wb.Worksheets.Item("Visium")

# The comment that currently sits on S1 (metadata_schema_id) needs to survive
# as the comment on Q1 once columns Q:R are removed and S shifts left. Update
# the Q1 comment's text in place (so authorship/shapeId bookkeeping is left
# alone) and drop the comments that won't survive the edit.
$finalQ1CommentText = $ws.Range("S1").Comment.Text()
$ws.Range("Q1").Comment.Text($finalQ1CommentText)
$ws.Range("R1").Comment.Delete()
$ws.Range("S1").Comment.Delete()

# Delete columns Q (preparation_instrument_vendor) and R
# (preparation_instrument_model) entirely; column S (metadata_schema_id)
# shifts left to become the new column Q, and the data validations that
# referenced the deleted lookup sheets go away with the columns.
$ws.Range("Q:R").EntireColumn.Delete()

# Restore the earlier pav:createdOn timestamp on the .metadata sheet.
$meta = $wb.Worksheets.Item(".metadata")
$meta.Range("C2").Value = "2023-11-01T15:37:30-07:00"

# Remove the now-unused lookup sheets.
$wb.Worksheets.Item("preparation_instrument_vendor").Delete()
$wb.Worksheets.Item("preparation_instrument_model").Delete()
